# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-mark cells whose new values are plain numeric-looking strings as Text
# so Excel stores them as text (matching the source data), not as numbers.
$ws.Range("D5,D6,D7,D9,D10,D11,D13,D15,D17,D19,D22,D23,D24,D25,D27,D28,D29,D30,D31,D32,D33,D35,D36,D37,D38,D39,D40,D42,D45,D46,D47,D48,D49,D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.555.04"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.339.67"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "304.93"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").Value = "101.37"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("D7").Value = "0.509"
$ws.Range("E7").Value = "  -3.35%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("D10").Value = "35.25"
$ws.Range("E10").Value = "  -2.71%  "
$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "6.81"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").Value = "2.703.46"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "15.65"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "2.383.29"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "0.809"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "43.452.43"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "11.85"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").Value = "68.34"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "238.01"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  -3.23%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "25.10"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").Value = "34.77"
$ws.Range("E28").Value = "  -4.75%  "
$ws.Range("D29").Value = "2.08"
$ws.Range("E29").Value = "  -5.64%  "
$ws.Range("D30").Value = "166.35"
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("D31").Value = "9.26"
$ws.Range("E31").Value = "  -3.59%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "5.07"
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("E34").Value = "  -4.73%  "
$ws.Range("D35").Value = "4.51"
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0707"
$ws.Range("E36").Value = "  -4.53%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "16.87"
$ws.Range("E37").Value = "  -7.88%  "
$ws.Range("D38").Value = "2.92"
$ws.Range("E38").Value = "  -7.03%  "
$ws.Range("D39").Value = "1.83"
$ws.Range("E39").Value = "  -6.43%  "
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").Value = "  -3.30%  "
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").Value = "2.41"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "1.986.40"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("D45").Value = "18.62"
$ws.Range("E45").Value = "  -8.91%  "
$ws.Range("D46").Value = "10.05"
$ws.Range("E46").Value = "  -3.96%  "
$ws.Range("D47").Value = "2.95"
$ws.Range("E47").Value = "  -7.80%  "
$ws.Range("D48").Value = "56.56"
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("D49").Value = "4.86"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("D50").Value = "2.565.46"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "1.55"
$ws.Range("E51").Value = "  -0.87%  "
